$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nowcast_R")

# Update the column header in Q1: "Wochentag" -> "Wochentag für RKI Tagesbericht"
$ws.Range("Q1").Value = "Wochentag für RKI Tagesbericht"

# The user selected A95:AL116 (covering the weekday-grouped statistics block in
# rows 102-114) and cleared its contents, removing the per-weekday averages /
# MAE stats / weekday-name labels while leaving the surrounding cell styling.
$ws.Activate()
$rng = $ws.Range("A95:AL116")
$rng.Select()
$rng.ClearContents()

# Recalculate so the dependent sheets (Auswertung Wochentage, MW + STD Schätzer)
# pick up the now-blank source cells.
$excel.Calculate()
